# Update essay count references from "70+" to "100+" and widen / reposition
# the related pill shapes on the client-deck "team" slide so the longer
# label still fits without overlapping its neighbour.

$p = $ppt.ActivePresentation

# --- Slide 6: report teaser paragraph -------------------------------------
$s6 = $p.Slides.Item(6)
$tb18 = $s6.Shapes.Item(18)
$tr = $tb18.TextFrame.TextRange
$tr.Text = $tr.Text.Replace("70+", "100+")

# --- Slide 8: bio paragraph + "70+ philosophical essays" pill -------------
$s8 = $p.Slides.Item(8)

# Bio paragraph ("Author of 70+ philosophical essays ...")
$bio = $s8.Shapes.Item(6)
$bioTr = $bio.TextFrame.TextRange
$bioTr.Text = $bioTr.Text.Replace("70+", "100+")

# Pill background rectangle behind "70+ philosophical essays" -> widen it
# (167.40001pt, not the exact 167.4pt, to land on the target EMU exactly
# once the host round-trips the value through a single-precision float)
$pillRect = $s8.Shapes.Item(15)
$pillRect.Width = 167.40001

# Pill label textbox "70+ philosophical essays" -> widen + retext
$pillText = $s8.Shapes.Item(16)
$pillText.Width = 153
$pillTextTr = $pillText.TextFrame.TextRange
$pillTextTr.Text = $pillTextTr.Text.Replace("70+", "100+")

# Next pill ("Off-grid validated") background rectangle shifts right to
# make room for the now-wider essay-count pill
$nextRect = $s8.Shapes.Item(17)
$nextRect.Left = 355.68001

# Next pill ("Off-grid validated") label textbox shifts right too
$nextText = $s8.Shapes.Item(18)
$nextText.Left = 362.88
